$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.70913233333333
$ws.Range("H2").Value = 32.127397
$ws.Range("I2").Value = 0.007451729107954897
$ws.Range("J2").Value = 0.007451729107954897
$ws.Range("M2").Value = 4.820639
$ws.Range("N2").Value = 14.461917
$ws.Range("O2").Value = 0.03139697817829975
$ws.Range("P2").Value = 0.03139697817829975
$ws.Range("Q2").Value = 51.62486098222767
$ws.Range("R2").Value = 464.623748840049
$ws.Range("S2").Value = 0.0002339617761930609
$ws.Range("T2").Value = 0.0002339617761930609
$ws.Range("G3").Value = 10.70913233333333
$ws.Range("H3").Value = 32.127397
$ws.Range("I3").Value = 0.007451729107954897
$ws.Range("J3").Value = 0.007451729107954897
$ws.Range("O3").Value = 0.2306102199252841
$ws.Range("P3").Value = 0.2306102199252841
$ws.Range("Q3").Value = 379.1836423593186
$ws.Range("R3").Value = 3412.652781233867
$ws.Range("S3").Value = 0.00171844488840912
$ws.Range("T3").Value = 0.00171844488840912
$ws.Range("G4").Value = 10.70913233333333
$ws.Range("H4").Value = 32.127397
$ws.Range("I4").Value = 0.007451729107954897
$ws.Range("J4").Value = 0.007451729107954897
$ws.Range("M4").Value = 51.32089766666667
$ws.Range("N4").Value = 153.962693
$ws.Range("O4").Value = 0.3342546712440172
$ws.Range("P4").Value = 0.3342546712440172
$ws.Range("Q4").Value = 549.6022845777912
$ws.Range("R4").Value = 4946.420561200121
$ws.Range("S4").Value = 0.002490775263178938
$ws.Range("T4").Value = 0.002490775263178938
$ws.Range("G5").Value = 10.70913233333333
$ws.Range("H5").Value = 32.127397
$ws.Range("I5").Value = 0.007451729107954897
$ws.Range("J5").Value = 0.007451729107954897
$ws.Range("M5").Value = 4.113383
$ws.Range("N5").Value = 12.340149
$ws.Range("O5").Value = 0.0267905969084159
$ws.Range("P5").Value = 0.0267905969084159
$ws.Range("Q5").Value = 44.05076288468367
$ws.Range("R5").Value = 396.456865962153
$ws.Range("S5").Value = 0.0001996362708019293
$ws.Range("T5").Value = 0.0001996362708019293
$ws.Range("G6").Value = 10.70913233333333
$ws.Range("H6").Value = 32.127397
$ws.Range("I6").Value = 0.007451729107954897
$ws.Range("J6").Value = 0.007451729107954897
$ws.Range("M6").Value = 57.87588766666666
$ws.Range("N6").Value = 173.627663
$ws.Range("O6").Value = 0.376947533743983
$ws.Range("P6").Value = 0.3769475337439831
$ws.Range("Q6").Value = 619.8005399314679
$ws.Range("R6").Value = 5578.204859383211
$ws.Range("S6").Value = 0.002808910909371849
$ws.Range("T6").Value = 0.002808910909371849
$ws.Range("I7").Value = 0.03290895798513831
$ws.Range("J7").Value = 0.03290895798513832
$ws.Range("M7").Value = 4.820639
$ws.Range("N7").Value = 14.461917
$ws.Range("O7").Value = 0.03139697817829975
$ws.Range("P7").Value = 0.03139697817829975
$ws.Range("Q7").Value = 227.9900888022217
$ws.Range("R7").Value = 2051.910799219995
$ws.Range("S7").Value = 0.001033241835729971
$ws.Range("T7").Value = 0.001033241835729971
$ws.Range("I8").Value = 0.03290895798513831
$ws.Range("J8").Value = 0.03290895798513832
$ws.Range("O8").Value = 0.2306102199252841
$ws.Range("P8").Value = 0.2306102199252841
$ws.Range("S8").Value = 0.00758914203846468
$ws.Range("T8").Value = 0.007589142038464683
$ws.Range("I9").Value = 0.03290895798513831
$ws.Range("J9").Value = 0.03290895798513832
$ws.Range("M9").Value = 51.32089766666667
$ws.Range("N9").Value = 153.962693
$ws.Range("O9").Value = 0.3342546712440172
$ws.Range("P9").Value = 0.3342546712440172
$ws.Range("Q9").Value = 2427.200214833151
$ws.Range("R9").Value = 21844.80193349835
$ws.Range("S9").Value = 0.01099997293230558
$ws.Range("T9").Value = 0.01099997293230558
$ws.Range("I10").Value = 0.03290895798513831
$ws.Range("J10").Value = 0.03290895798513832
$ws.Range("M10").Value = 4.113383
$ws.Range("N10").Value = 12.340149
$ws.Range("O10").Value = 0.0267905969084159
$ws.Range("P10").Value = 0.02679059690841591
$ws.Range("Q10").Value = 194.5407145085017
$ws.Range("R10").Value = 1750.866430576515
$ws.Range("S10").Value = 0.0008816506280558354
$ws.Range("T10").Value = 0.0008816506280558358
$ws.Range("I11").Value = 0.03290895798513831
$ws.Range("J11").Value = 0.03290895798513832
$ws.Range("M11").Value = 57.87588766666666
$ws.Range("N11").Value = 173.627663
$ws.Range("O11").Value = 0.376947533743983
$ws.Range("P11").Value = 0.3769475337439831
$ws.Range("Q11").Value = 2737.215702862367
$ws.Range("R11").Value = 24634.9413257613
$ws.Range("S11").Value = 0.01240495055058224
$ws.Range("T11").Value = 0.01240495055058225
$ws.Range("G12").Value = 411.37678
$ws.Range("H12").Value = 1234.13034
$ws.Range("I12").Value = 0.2862480573072345
$ws.Range("J12").Value = 0.2862480573072345
$ws.Range("M12").Value = 4.820639
$ws.Range("N12").Value = 14.461917
$ws.Range("O12").Value = 0.03139697817829975
$ws.Range("P12").Value = 0.03139697817829975
$ws.Range("Q12").Value = 1983.09894936242
$ws.Range("R12").Value = 17847.89054426178
$ws.Range("S12").Value = 0.008987324008855936
$ws.Range("T12").Value = 0.008987324008855936
$ws.Range("G13").Value = 411.37678
$ws.Range("H13").Value = 1234.13034
$ws.Range("I13").Value = 0.2862480573072345
$ws.Range("J13").Value = 0.2862480573072345
$ws.Range("O13").Value = 0.2306102199252841
$ws.Range("P13").Value = 0.2306102199252841
$ws.Range("Q13").Value = 14565.82484623153
$ws.Range("R13").Value = 131092.4236160837
$ws.Range("S13").Value = 0.06601172744880666
$ws.Range("T13").Value = 0.06601172744880666
$ws.Range("G14").Value = 411.37678
$ws.Range("H14").Value = 1234.13034
$ws.Range("I14").Value = 0.2862480573072345
$ws.Range("J14").Value = 0.2862480573072345
$ws.Range("M14").Value = 51.32089766666667
$ws.Range("N14").Value = 153.962693
$ws.Range("O14").Value = 0.3342546712440172
$ws.Range("P14").Value = 0.3342546712440172
$ws.Range("Q14").Value = 21112.22562882285
$ws.Range("R14").Value = 190010.0306594056
$ws.Range("S14").Value = 0.09567975028946825
$ws.Range("T14").Value = 0.09567975028946825
$ws.Range("G15").Value = 411.37678
$ws.Range("H15").Value = 1234.13034
$ws.Range("I15").Value = 0.2862480573072345
$ws.Range("J15").Value = 0.2862480573072345
$ws.Range("M15").Value = 4.113383
$ws.Range("N15").Value = 12.340149
$ws.Range("O15").Value = 0.0267905969084159
$ws.Range("P15").Value = 0.02679059690841591
$ws.Range("Q15").Value = 1692.15025344674
$ws.Range("R15").Value = 15229.35228102066
$ws.Range("S15").Value = 0.007668756319135254
$ws.Range("T15").Value = 0.007668756319135256
$ws.Range("G16").Value = 411.37678
$ws.Range("H16").Value = 1234.13034
$ws.Range("I16").Value = 0.2862480573072345
$ws.Range("J16").Value = 0.2862480573072345
$ws.Range("M16").Value = 57.87588766666666
$ws.Range("N16").Value = 173.627663
$ws.Range("O16").Value = 0.376947533743983
$ws.Range("P16").Value = 0.3769475337439831
$ws.Range("Q16").Value = 23808.79630795505
$ws.Range("R16").Value = 214279.1667715954
$ws.Range("S16").Value = 0.1079004992409684
$ws.Range("T16").Value = 0.1079004992409684
$ws.Range("G17").Value = 173.2560603333334
$ws.Range("H17").Value = 519.768181
$ws.Range("I17").Value = 0.12055666021578
$ws.Range("J17").Value = 0.12055666021578
$ws.Range("M17").Value = 4.820639
$ws.Range("N17").Value = 14.461917
$ws.Range("O17").Value = 0.03139697817829975
$ws.Range("P17").Value = 0.03139697817829975
$ws.Range("Q17").Value = 835.2049214292198
$ws.Range("R17").Value = 7516.844292862977
$ws.Range("S17").Value = 0.003785114830043541
$ws.Range("T17").Value = 0.003785114830043542
$ws.Range("G18").Value = 173.2560603333334
$ws.Range("H18").Value = 519.768181
$ws.Range("I18").Value = 0.12055666021578
$ws.Range("J18").Value = 0.12055666021578
$ws.Range("O18").Value = 0.2306102199252841
$ws.Range("P18").Value = 0.2306102199252841
$ws.Range("Q18").Value = 6134.564591524721
$ws.Range("R18").Value = 55211.08132372249
$ws.Range("S18").Value = 0.02780159792581877
$ws.Range("T18").Value = 0.02780159792581877
$ws.Range("G19").Value = 173.2560603333334
$ws.Range("H19").Value = 519.768181
$ws.Range("I19").Value = 0.12055666021578
$ws.Range("J19").Value = 0.12055666021578
$ws.Range("M19").Value = 51.32089766666667
$ws.Range("N19").Value = 153.962693
$ws.Range("O19").Value = 0.3342546712440172
$ws.Range("P19").Value = 0.3342546712440172
$ws.Range("Q19").Value = 8891.656542496827
$ws.Range("R19").Value = 80024.90888247144
$ws.Range("S19").Value = 0.04029662682670222
$ws.Range("T19").Value = 0.04029662682670223
$ws.Range("G20").Value = 173.2560603333334
$ws.Range("H20").Value = 519.768181
$ws.Range("I20").Value = 0.12055666021578
$ws.Range("J20").Value = 0.12055666021578
$ws.Range("M20").Value = 4.113383
$ws.Range("N20").Value = 12.340149
$ws.Range("O20").Value = 0.0267905969084159
$ws.Range("P20").Value = 0.02679059690841591
$ws.Range("Q20").Value = 712.6685332221077
$ws.Range("R20").Value = 6414.01679899897
$ws.Range("S20").Value = 0.003229784888465822
$ws.Range("T20").Value = 0.003229784888465823
$ws.Range("G21").Value = 173.2560603333334
$ws.Range("H21").Value = 519.768181
$ws.Range("I21").Value = 0.12055666021578
$ws.Range("J21").Value = 0.12055666021578
$ws.Range("M21").Value = 57.87588766666666
$ws.Range("N21").Value = 173.627663
$ws.Range("O21").Value = 0.376947533743983
$ws.Range("P21").Value = 0.3769475337439831
$ws.Range("Q21").Value = 10027.34828542122
$ws.Range("R21").Value = 90246.13456879099
$ws.Range("S21").Value = 0.04544353574474962
$ws.Range("T21").Value = 0.04544353574474963
$ws.Range("G22").Value = 794.4973246666667
$ws.Range("H22").Value = 2383.491974
$ws.Range("I22").Value = 0.5528345953838922
$ws.Range("J22").Value = 0.5528345953838923
$ws.Range("M22").Value = 4.820639
$ws.Range("N22").Value = 14.461917
$ws.Range("O22").Value = 0.03139697817829975
$ws.Range("P22").Value = 0.03139697817829975
$ws.Range("Q22").Value = 3829.984788683796
$ws.Range("R22").Value = 34469.86309815416
$ws.Range("S22").Value = 0.01735733572747724
$ws.Range("T22").Value = 0.01735733572747724
$ws.Range("G23").Value = 794.4973246666667
$ws.Range("H23").Value = 2383.491974
$ws.Range("I23").Value = 0.5528345953838922
$ws.Range("J23").Value = 0.5528345953838923
$ws.Range("O23").Value = 0.2306102199252841
$ws.Range("P23").Value = 0.2306102199252841
$ws.Range("Q23").Value = 28131.16693629185
$ws.Range("R23").Value = 253180.5024266267
$ws.Range("S23").Value = 0.1274893076237848
$ws.Range("T23").Value = 0.1274893076237849
$ws.Range("G24").Value = 794.4973246666667
$ws.Range("H24").Value = 2383.491974
$ws.Range("I24").Value = 0.5528345953838922
$ws.Range("J24").Value = 0.5528345953838923
$ws.Range("M24").Value = 51.32089766666667
$ws.Range("N24").Value = 153.962693
$ws.Range("O24").Value = 0.3342546712440172
$ws.Range("P24").Value = 0.3342546712440172
$ws.Range("Q24").Value = 40774.31589565844
$ws.Range("R24").Value = 366968.843060926
$ws.Range("S24").Value = 0.1847875459323622
$ws.Range("T24").Value = 0.1847875459323622
$ws.Range("G25").Value = 794.4973246666667
$ws.Range("H25").Value = 2383.491974
$ws.Range("I25").Value = 0.5528345953838922
$ws.Range("J25").Value = 0.5528345953838923
$ws.Range("M25").Value = 4.113383
$ws.Range("N25").Value = 12.340149
$ws.Range("O25").Value = 0.0267905969084159
$ws.Range("P25").Value = 0.02679059690841591
$ws.Range("Q25").Value = 3268.071788829347
$ws.Range("R25").Value = 29412.64609946413
$ws.Range("S25").Value = 0.01481076880195706
$ws.Range("T25").Value = 0.01481076880195707
$ws.Range("G26").Value = 794.4973246666667
$ws.Range("H26").Value = 2383.491974
$ws.Range("I26").Value = 0.5528345953838922
$ws.Range("J26").Value = 0.5528345953838923
$ws.Range("M26").Value = 57.87588766666666
$ws.Range("N26").Value = 173.627663
$ws.Range("O26").Value = 0.376947533743983
$ws.Range("P26").Value = 0.3769475337439831
$ws.Range("Q26").Value = 45982.23791387519
$ws.Range("R26").Value = 413840.1412248767
$ws.Range("S26").Value = 0.2083896372983109
$ws.Range("T26").Value = 0.208389637298311

Write-Output "Applied 278 cell updates"